$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 16: Using Your Arcane Powers for Fun and Profit / Ash Picatrix
$ws.Range("H16").Value = 5500
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -770

# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 596035.9
$ws.Range("J17").Value = 596035.9
$ws.Range("L17").Value = 1788107.7
$ws.Range("N17").Value = -1788443.7

# Row 21: Book and a Hard Place / Engraved Hard Leather Grimoire
$ws.Range("H21").Value = 21254.5
$ws.Range("J21").Value = 21254.5
$ws.Range("L21").Value = 21254.5
$ws.Range("N21").Value = -22190.5

# Row 23: There's Something about Bury / Hard Leather Grimoire
$ws.Range("H23").Value = 21254.5
$ws.Range("J23").Value = 21254.5
$ws.Range("L23").Value = 21254.5
$ws.Range("N23").Value = -21722.5

# Row 32: Automata for the People / Crab Oil
$ws.Range("H32").Value = 2498.5
$ws.Range("J32").Value = 2498
$ws.Range("L32").Value = 2498
$ws.Range("N32").Value = -3150

# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 3243.375
$ws.Range("J40").Value = 3839.8
$ws.Range("L40").Value = 3839.8
$ws.Range("N40").Value = -4189.8

# Row 51: A Bile Business / Shark Oil
$ws.Range("H51").Value = 83333336
$ws.Range("J51").Value = 83333336
$ws.Range("L51").Value = 83333336
$ws.Range("N51").Value = -83334304

# Row 97: Materia Worth / Potent Spiritbond Potion
$ws.Range("H97").Value = 1801.1666
$ws.Range("J97").Value = 1801.1666
$ws.Range("L97").Value = 5403.4998
$ws.Range("N97").Value = -6395.4998

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 571142.7
$ws.Range("I132").Value = 659426.25
$ws.Range("K132").Value = 1978278.75
$ws.Range("M132").Value = -1975748.75

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 1392.2778
$ws.Range("I141").Value = 1392.2778
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4176.8334
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 1003.1666
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 1823.2354
$ws.Range("I2").Value = 1728.3125
$ws.Range("K2").Value = 1728.3125
$ws.Range("M2").Value = -1615.3125

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 4670.7075
$ws.Range("I32").Value = 3944.2295
$ws.Range("J32").Value = 15749.5
$ws.Range("K32").Value = 3944.2295
$ws.Range("L32").Value = 15749.5
$ws.Range("M32").Value = -3657.2295
$ws.Range("N32").Value = -16323.5

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 1718.5555
$ws.Range("I61").Value = 1718.5555
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1718.5555
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1506.5555
$ws.Range("N61").Value = $null

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1173.8914
$ws.Range("I74").Value = 962.4872
$ws.Range("J74").Value = 2351.7144
$ws.Range("K74").Value = 962.4872
$ws.Range("L74").Value = 2351.7144
$ws.Range("M74").Value = -88.48720000000003
$ws.Range("N74").Value = -4099.7144

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1173.8914
$ws.Range("I77").Value = 962.4872
$ws.Range("J77").Value = 2351.7144
$ws.Range("K77").Value = 4812.436
$ws.Range("L77").Value = 11758.572
$ws.Range("M77").Value = -444.4359999999997
$ws.Range("N77").Value = -20494.572

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 1823.2354
$ws.Range("I116").Value = 1728.3125
$ws.Range("K116").Value = 1728.3125
$ws.Range("M116").Value = 565.6875

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 5452.65
$ws.Range("I132").Value = 4852.0566
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 14556.1698
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -12026.1698
$ws.Range("N132").Value = -35060

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1718.5555
$ws.Range("I136").Value = 1718.5555
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5155.666499999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2605.666499999999
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("BSM")
# Row 2: Proly Hatchet / Bronze Hatchet
$ws.Range("H2").Value = 39993.5
$ws.Range("J2").Value = 39993.5
$ws.Range("L2").Value = 39993.5
$ws.Range("N2").Value = -40219.5

# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 1823.2354
$ws.Range("I3").Value = 1728.3125
$ws.Range("K3").Value = 1728.3125
$ws.Range("M3").Value = -1614.3125

# Row 6: The Unkindest Cut / Bronze Saw
$ws.Range("H6").Value = 28028
$ws.Range("J6").Value = 28028
$ws.Range("L6").Value = 28028
$ws.Range("N6").Value = -28254

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 2749.1924
$ws.Range("I86").Value = 2746.8572
$ws.Range("J86").Value = 2759
$ws.Range("K86").Value = 2746.8572
$ws.Range("L86").Value = 2759
$ws.Range("M86").Value = -1623.8572
$ws.Range("N86").Value = -5005

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 2749.1924
$ws.Range("I89").Value = 2746.8572
$ws.Range("J89").Value = 2759
$ws.Range("K89").Value = 13734.286
$ws.Range("L89").Value = 13795
$ws.Range("M89").Value = -8118.286
$ws.Range("N89").Value = -25027

# Row 124: History of the Hrothgar / High Durium Bayonet
$ws.Range("H124").Value = 118332.664
$ws.Range("J124").Value = 118332.664
$ws.Range("L124").Value = 118332.664
$ws.Range("N124").Value = -128152.664

# Row 125: Archon of His Eye / High Durium Knives
$ws.Range("H125").Value = 45000
$ws.Range("J125").Value = 45000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840

# Row 130: Annals of the Empire I / Chondrite Magitek Axe
$ws.Range("H130").Value = 80000
$ws.Range("J130").Value = 80000
$ws.Range("L130").Value = 80000
$ws.Range("N130").Value = -90040

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1248.2413
$ws.Range("I134").Value = 1081.4445
$ws.Range("K134").Value = 3244.3335
$ws.Range("M134").Value = -709.3335000000002

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 1995.5454
$ws.Range("I58").Value = 1923.4286
$ws.Range("K58").Value = 1923.4286
$ws.Range("M58").Value = -1720.4286

# Row 86: Birch, Please / Birch Lumber
$ws.Range("H86").Value = 3238
$ws.Range("I86").Value = 1902.8334
$ws.Range("K86").Value = 1902.8334
$ws.Range("M86").Value = -779.8334

# Row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws.Range("H89").Value = 3238
$ws.Range("I89").Value = 1902.8334
$ws.Range("K89").Value = 9514.166999999999
$ws.Range("M89").Value = -3898.166999999999

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 41339.19
$ws.Range("I122").Value = 69323.60000000001
$ws.Range("K122").Value = 207970.8
$ws.Range("M122").Value = -205520.8

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1917.55
$ws.Range("I134").Value = 1715.8422
$ws.Range("K134").Value = 5147.5266
$ws.Range("M134").Value = -2612.5266

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 1995.5454
$ws.Range("I136").Value = 1923.4286
$ws.Range("K136").Value = 5770.2858
$ws.Range("M136").Value = -3220.2858

$ws = $wb.Worksheets.Item("CUL")
# Row 114: One Last Meal / Mushroom Saute
$ws.Range("H114").Value = 22223844
$ws.Range("I114").Value = 22223844
$ws.Range("K114").Value = 66671532
$ws.Range("M114").Value = -66668278

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 58199.11
$ws.Range("J131").Value = 3999.5
$ws.Range("L131").Value = 11998.5
$ws.Range("N131").Value = -22078.5

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers / Copper Ingot
$ws.Range("H2").Value = 276.23077
$ws.Range("I2").Value = 375.2
$ws.Range("J2").Value = 172.05263
$ws.Range("K2").Value = 375.2
$ws.Range("L2").Value = 172.05263
$ws.Range("M2").Value = -262.2
$ws.Range("N2").Value = -398.05263

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 1588.3226
$ws.Range("I102").Value = 1646.2069
$ws.Range("K102").Value = 1646.2069
$ws.Range("M102").Value = -24.20689999999991

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 2434.257
$ws.Range("I132").Value = 2434.257
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7302.771000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4772.771000000001
$ws.Range("N132").Value = $null

# Row 136: Shiny and Good / Pink Beryl
$ws.Range("H136").Value = 34636.395
$ws.Range("J136").Value = 34636.395
$ws.Range("L136").Value = 103909.185
$ws.Range("N136").Value = -109009.185

$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 168682.08
$ws.Range("I61").Value = 114057.11
$ws.Range("J61").Value = 332557
$ws.Range("K61").Value = 114057.11
$ws.Range("L61").Value = 332557
$ws.Range("M61").Value = -113855.11
$ws.Range("N61").Value = -332961

# Row 64: Glorified Hole-punchers / Archaeoskin Gloves of Aiming
$ws.Range("H64").Value = 46935.4
$ws.Range("J64").Value = 46935.4
$ws.Range("L64").Value = 46935.4
$ws.Range("N64").Value = -47385.4

# Row 67: Treat Them with Kid Gloves (L) / Archaeoskin Gloves of Aiming
$ws.Range("H67").Value = 46935.4
$ws.Range("J67").Value = 46935.4
$ws.Range("L67").Value = 46935.4
$ws.Range("N67").Value = -48495.4

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 168682.08
$ws.Range("I113").Value = 114057.11
$ws.Range("J113").Value = 332557
$ws.Range("K113").Value = 114057.11
$ws.Range("L113").Value = 332557
$ws.Range("M113").Value = -111887.11
$ws.Range("N113").Value = -336897

# Row 131: For What Was Gleaned / Ophiotauroskin Wristband of Gathering
$ws.Range("H131").Value = 85000
$ws.Range("J131").Value = 85000
$ws.Range("L131").Value = 85000
$ws.Range("N131").Value = -95080

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 3058.7415
$ws.Range("J132").Value = 4188.7144
$ws.Range("L132").Value = 12566.1432
$ws.Range("N132").Value = -17626.1432

$ws = $wb.Worksheets.Item("WVR")
# Row 124: Hot Heads / Almasty Serge Hat of Casting
$ws.Range("H124").Value = 88379
$ws.Range("J124").Value = 88379
$ws.Range("L124").Value = 88379
$ws.Range("N124").Value = -98199

# Row 128: Lightening Up / Scarlet Moko Gaskins of the Rising Dragon
$ws.Range("H128").Value = 69800
$ws.Range("J128").Value = 69800
$ws.Range("L128").Value = 69800
$ws.Range("N128").Value = -79760

# Row 129: Lifetime of Gleaning / Scarlet Moko Beret of Gathering
$ws.Range("H129").Value = 98398.60000000001
$ws.Range("J129").Value = 98398.60000000001
$ws.Range("L129").Value = 98398.60000000001
$ws.Range("N129").Value = -108398.6

# Row 130: Skill Cap / AR-Caean Velvet Cap of Maiming
$ws.Range("H130").Value = 42427.855
$ws.Range("J130").Value = 42427.855
$ws.Range("L130").Value = 42427.855
$ws.Range("N130").Value = -52467.855

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 3133.8655
$ws.Range("I132").Value = 2545.9556
$ws.Range("K132").Value = 7637.8668
$ws.Range("M132").Value = -5107.8668
